# Update ticket-count / min-price figures on the "展览" and "全部类型" sheets.
# (The "演出" and "本地生活" sheets are untouched by this revision.)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 80
$ws1.Range("G2").Value = "不可售"

$ws1.Range("G3").Value = 80

$ws1.Range("F4").Value = 3013
$ws1.Range("F5").Value = 459
$ws1.Range("F7").Value = 25
$ws1.Range("F9").Value = 5
$ws1.Range("F10").Value = 14371
$ws1.Range("F11").Value = 151
$ws1.Range("F12").Value = 118
$ws1.Range("F13").Value = 5748
$ws1.Range("F19").Value = 11
$ws1.Range("F25").Value = 10537
$ws1.Range("F27").Value = 56
$ws1.Range("F28").Value = 78

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 80
$ws4.Range("G2").Value = "不可售"

$ws4.Range("G3").Value = 80

$ws4.Range("F5").Value = 3013
$ws4.Range("F6").Value = 459
$ws4.Range("F8").Value = 25
$ws4.Range("F10").Value = 5
$ws4.Range("F11").Value = 14371
$ws4.Range("F12").Value = 151
$ws4.Range("F13").Value = 118
$ws4.Range("F14").Value = 5748
$ws4.Range("F20").Value = 11
$ws4.Range("F27").Value = 10537
$ws4.Range("F29").Value = 56
$ws4.Range("F30").Value = 78
